# Schematic and PCB V1 done
# - Rename the "1x2 conn" part to "1x2 conn female" and correct its quantity.
# - Add a new BOM line for the matching "1x2 conn male" connector.
# - Move the active selection to D18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "1x2 conn" -> "1x2 conn female", quantity 3 -> 2 -------------
$ws.Range("E8").Value = "1x2 conn female"
$ws.Range("C8").Value = 2

# --- Row 13 (new line): male counterpart of the 1x2 connector ------------
$a13 = $ws.Range("A13")
$a13.Value = "649-1012938090202ALF`n"
$a13.Font.Size = 10
$a13.Font.Bold = $false
$a13.WrapText = $true
$a13.HorizontalAlignment = -4131
$a13.VerticalAlignment = -4160

$b13 = $ws.Range("B13")
$b13.Value = "10129380-902002ALF`n`n"
$b13.Font.Size = 10
$b13.Font.Bold = $false
$b13.WrapText = $true
$b13.HorizontalAlignment = -4131
$b13.VerticalAlignment = -4160
$b13.Characters(1, 19).Font.Size = 10
$b13.Characters(1, 19).Font.Bold = $false
$b13.Characters(20, 1).Font.Size = 16
$b13.Characters(20, 1).Font.Bold = $false

$ws.Range("C13").Value = 2

$f13 = $ws.Range("F13")
$f13.Value = "1x2 conn male"
$f13.Font.Size = 10
$f13.Font.Bold = $false
$f13.WrapText = $false
$f13.HorizontalAlignment = 1
$f13.VerticalAlignment = -4107

# Keep row 13's height consistent with the rest of the sheet (wrapping must
# not trigger an auto row-height change).
$ws.Rows.Item(13).RowHeight = 12.8

# --- Move the active selection, matching where the author left off -------
$ws.Range("D18").Select() | Out-Null
